$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of price values look fully numeric but carry a significant
# trailing zero (e.g. 0.670, 98.10) that Excel's normal numeric-literal
# inference would silently drop (-> 0.67 / 98.1). Pre-format just those
# cells as Text so the literal digits survive exactly as scraped.
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'

$ws.Range('D2').Value = '37.171.13'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '2.057.66'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '250.31'
$ws.Range('E5').Value = '  +0.62%  '
$ws.Range('D6').Value = '0.670'
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('D7').Value = '61.09'
$ws.Range('E7').Value = '  +10.32%  '
$ws.Range('D9').Value = '0.386'
$ws.Range('E9').Value = '  +2.02%  '
$ws.Range('D10').Value = '0.0792'
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('E11').Value = '  +1.96%  '
$ws.Range('D12').Value = '16.28'
$ws.Range('E12').Value = '  +7.73%  '
$ws.Range('D13').Value = '2.355.75'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '0.826'
$ws.Range('E14').Value = '  +1.29%  '
$ws.Range('D15').Value = '5.72'
$ws.Range('E15').Value = '  +9.16%  '
$ws.Range('D16').Value = '2.061.64'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').Value = '18.04'
$ws.Range('E17').Value = '  +27.05%  '
$ws.Range('D18').Value = '37.152.29'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = '75.36'
$ws.Range('E19').Value = '  +4.10%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -3.33%  '
$ws.Range('D21').Value = '5.44'
$ws.Range('E21').Value = '  +1.07%  '
$ws.Range('D22').Value = '239.73'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  -1.04%  '
$ws.Range('D25').Value = '2.22'
$ws.Range('E25').Value = '  +12.70%  '
$ws.Range('D26').Value = '169.44'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').Value = '9.43'
$ws.Range('E27').Value = '  +4.25%  '
$ws.Range('D28').Value = '20.06'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('E29').Value = '  +1.84%  '
$ws.Range('E30').Value = '  +9.52%  '
$ws.Range('D31').Value = '4.83'
$ws.Range('E31').Value = '  +5.88%  '
$ws.Range('D32').Value = '0.0621'
$ws.Range('E32').Value = '  -0.53%  '
$ws.Range('D33').Value = '4.56'
$ws.Range('E33').Value = '  +4.36%  '
$ws.Range('D34').Value = '0.0898'
$ws.Range('E34').Value = '  +4.94%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').Value = '1.74'
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('D38').Value = '0.109'
$ws.Range('E38').Value = '  +3.64%  '
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('D40').Value = '5.33'
$ws.Range('E40').Value = '  +32.48%  '
$ws.Range('E41').Value = '  +14.29%  '
$ws.Range('D42').Value = '18.22'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').Value = '0.0224'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '98.10'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '2.49'
$ws.Range('E46').Value = '  +3.23%  '
$ws.Range('D47').Value = '1.296.81'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').Value = '6.87'
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('D50').Value = '2.240.83'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = '3.57'
$ws.Range('E51').Value = '  -15.71%  '
